$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear header detail cells AY2 and AY3 (duplicate tax-currency detail columns removed)
$ws.Range("AY2").ClearContents()
$ws.Range("AY3").ClearContents()

# Per-row updates: move currency from AA to Y, clear AA, set AO tax-rate text, clear AY file name
$rows = @(
    @{Row=5; Currency="AUD"; Rate="10.00 %"}
    @{Row=6; Currency="AUD"; Rate="10.00 %"}
    @{Row=7; Currency="AUD"; Rate="10.00 %"}
    @{Row=8; Currency="AUD"; Rate="10.00 %"}
    @{Row=9; Currency="AUD"; Rate="10.00 %"}
    @{Row=10; Currency="AUD"; Rate="10.00 %"}
    @{Row=11; Currency="AUD"; Rate="10.00 %"}
    @{Row=12; Currency="AUD"; Rate="10.00 %"}
    @{Row=13; Currency="AUD"; Rate="10.00 %"}
    @{Row=14; Currency="AUD"; Rate="10.00 %"}
    @{Row=15; Currency="PLN"; Rate="0%"}
    @{Row=16; Currency="DKK"; Rate="0%"}
    @{Row=17; Currency="PLN"; Rate="0%"}
    @{Row=18; Currency="EUR"; Rate="0%"}
    @{Row=19; Currency="RON"; Rate="0%"}
    @{Row=20; Currency="SEK"; Rate="0%"}
    @{Row=21; Currency="CHF"; Rate="0%"}
    @{Row=22; Currency="CZK"; Rate="0%"}
    @{Row=23; Currency="CZK"; Rate="0%"}
    @{Row=24; Currency="SEK"; Rate="0%"}
    @{Row=25; Currency="GBP"; Rate="20.00 %"}
    @{Row=26; Currency="GBP"; Rate="20.00 %"}
    @{Row=27; Currency="GBP"; Rate="20.00 %"}
    @{Row=28; Currency="GBP"; Rate="20.00 %"}
    @{Row=29; Currency="GBP"; Rate="20.00 %"}
    @{Row=30; Currency="GBP"; Rate="20.00 %"}
    @{Row=31; Currency="GBP"; Rate="20.00 %"}
    @{Row=32; Currency="GBP"; Rate="20.00 %"}
    @{Row=33; Currency="GBP"; Rate="20.00 %"}
    @{Row=34; Currency="GBP"; Rate="20.00 %"}
    @{Row=35; Currency="PLN"; Rate="0%"}
    @{Row=36; Currency="EUR"; Rate="0%"}
    @{Row=37; Currency="NOK"; Rate="0%"}
    @{Row=38; Currency="HUF"; Rate="0%"}
    @{Row=39; Currency="HUF"; Rate="0%"}
    @{Row=40; Currency="DKK"; Rate="0%"}
    @{Row=41; Currency="SEK"; Rate="0%"}
    @{Row=42; Currency="RON"; Rate="0%"}
    @{Row=43; Currency="RON"; Rate="0%"}
    @{Row=44; Currency="EUR"; Rate="0%"}
    @{Row=45; Currency="GBP"; Rate="20.00 %"}
    @{Row=46; Currency="GBP"; Rate="20.00 %"}
    @{Row=47; Currency="GBP"; Rate="20.00 %"}
    @{Row=48; Currency="GBP"; Rate="20.00 %"}
    @{Row=49; Currency="GBP"; Rate="20.00 %"}
    @{Row=50; Currency="GBP"; Rate="20.00 %"}
    @{Row=51; Currency="GBP"; Rate="20.00 %"}
    @{Row=52; Currency="GBP"; Rate="20.00 %"}
    @{Row=53; Currency="GBP"; Rate="20.00 %"}
    @{Row=54; Currency="GBP"; Rate="20.00 %"}
    @{Row=55; Currency="USD"; Rate="0%"}
    @{Row=56; Currency="USD"; Rate="0%"}
    @{Row=57; Currency="USD"; Rate="0%"}
    @{Row=58; Currency="USD"; Rate="0%"}
    @{Row=59; Currency="USD"; Rate="0%"}
    @{Row=60; Currency="USD"; Rate="0%"}
    @{Row=61; Currency="USD"; Rate="0%"}
    @{Row=62; Currency="USD"; Rate="0%"}
    @{Row=63; Currency="USD"; Rate="0%"}
    @{Row=64; Currency="USD"; Rate="0%"}
    @{Row=65; Currency="USD"; Rate="0%"}
    @{Row=66; Currency="USD"; Rate="0%"}
    @{Row=67; Currency="USD"; Rate="0%"}
    @{Row=68; Currency="USD"; Rate="0%"}
    @{Row=69; Currency="USD"; Rate="0%"}
    @{Row=70; Currency="USD"; Rate="0%"}
    @{Row=71; Currency="USD"; Rate="0%"}
    @{Row=72; Currency="USD"; Rate="0%"}
    @{Row=73; Currency="USD"; Rate="0%"}
    @{Row=74; Currency="USD"; Rate="0%"}
    @{Row=75; Currency="CAD"; Rate="0%"}
    @{Row=76; Currency="CAD"; Rate="0%"}
    @{Row=77; Currency="CAD"; Rate="0%"}
    @{Row=78; Currency="CAD"; Rate="0%"}
    @{Row=79; Currency="CAD"; Rate="0%"}
    @{Row=80; Currency="CAD"; Rate="0%"}
    @{Row=81; Currency="CAD"; Rate="0%"}
    @{Row=82; Currency="CAD"; Rate="0%"}
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 25).Value = $item.Currency   # column Y
    $ws.Cells.Item($r, 27).ClearContents()          # column AA
    $aoCell = $ws.Cells.Item($r, 41)                # column AO
    $aoCell.NumberFormat = "@"
    $aoCell.Value = $item.Rate
    $ws.Cells.Item($r, 51).ClearContents()          # column AY
}
